# Fruta / hortaliza, semanal
# Adds the latest weekly price observations (new rows 246-247, date 45013,
# "Región del Maule") for "Sandia" at Terminal Hortofrutícola Agro Chillán,
# pushing the previously-existing rows down by two positions so the new
# week's data sits at the top of this block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 248, shifting the
# old rows 248:269 down to 250:271.
$ws.Range("A248:R249").EntireRow.Insert()

# Columns that stay constant across this whole "Sandia" block.
$mercadoId = 7
$mercado = "Terminal Hortofrutícola Agro Chillán"
$region = "Ñuble"
$codreg = 16
$categoriaId = 100112028
$categoria = "Sandia"
$variedad = "Sin especificar"
$unidad = "$/unidad"
$kgOUnidades = 1
$clasificacion = "Hortaliza"

# --- New row 246: "Primera" quality, week of 2023-03-28 (serial 45013) ---
$ws.Range("A246").Value = $mercadoId
$ws.Range("B246").Value = $mercado
$ws.Range("C246").Value = $region
$ws.Range("D246").Value = 45013
$ws.Range("E246").Value = $codreg
$ws.Range("F246").Value = $categoriaId
$ws.Range("G246").Value = $categoria
$ws.Range("H246").Value = $variedad
$ws.Range("I246").Value = "Primera"
$ws.Range("J246").Value = 100
$ws.Range("K246").Value = 2000
$ws.Range("L246").Value = 2000
$ws.Range("M246").Value = 2000
$ws.Range("N246").Value = $unidad
$ws.Range("O246").Value = "Región del Maule"
$ws.Range("P246").Value = 2000
$ws.Range("Q246").Value = $kgOUnidades
$ws.Range("R246").Value = $clasificacion

# --- New row 247: "Segunda" quality, week of 2023-03-28 (serial 45013) ---
$ws.Range("A247").Value = $mercadoId
$ws.Range("B247").Value = $mercado
$ws.Range("C247").Value = $region
$ws.Range("D247").Value = 45013
$ws.Range("E247").Value = $codreg
$ws.Range("F247").Value = $categoriaId
$ws.Range("G247").Value = $categoria
$ws.Range("H247").Value = $variedad
$ws.Range("I247").Value = "Segunda"
$ws.Range("J247").Value = 100
$ws.Range("K247").Value = 1500
$ws.Range("L247").Value = 1500
$ws.Range("M247").Value = 1500
$ws.Range("N247").Value = $unidad
$ws.Range("O247").Value = "Región del Maule"
$ws.Range("P247").Value = 1500
$ws.Range("Q247").Value = $kgOUnidades
$ws.Range("R247").Value = $clasificacion

# --- Row 248 now holds what used to be row 246: "Extra" quality, 44610 ---
$ws.Range("A248").Value = $mercadoId
$ws.Range("B248").Value = $mercado
$ws.Range("C248").Value = $region
$ws.Range("D248").Value = 44610
$ws.Range("E248").Value = $codreg
$ws.Range("F248").Value = $categoriaId
$ws.Range("G248").Value = $categoria
$ws.Range("H248").Value = $variedad
$ws.Range("I248").Value = "Extra"
$ws.Range("J248").Value = 500
$ws.Range("K248").Value = 2500
$ws.Range("L248").Value = 2500
$ws.Range("M248").Value = 2500
$ws.Range("N248").Value = $unidad
$ws.Range("O248").Value = "Región de O'Higgins"
$ws.Range("P248").Value = 2500
$ws.Range("Q248").Value = $kgOUnidades
$ws.Range("R248").Value = $clasificacion

# --- Row 249 now holds what used to be row 247: "Primera" quality, 44610 ---
$ws.Range("A249").Value = $mercadoId
$ws.Range("B249").Value = $mercado
$ws.Range("C249").Value = $region
$ws.Range("D249").Value = 44610
$ws.Range("E249").Value = $codreg
$ws.Range("F249").Value = $categoriaId
$ws.Range("G249").Value = $categoria
$ws.Range("H249").Value = $variedad
$ws.Range("I249").Value = "Primera"
$ws.Range("J249").Value = 600
$ws.Range("K249").Value = 2000
$ws.Range("L249").Value = 2200
$ws.Range("M249").Value = 2100
$ws.Range("N249").Value = $unidad
$ws.Range("O249").Value = "Región de O'Higgins"
$ws.Range("P249").Value = 2100
$ws.Range("Q249").Value = $kgOUnidades
$ws.Range("R249").Value = $clasificacion

# Apply the same date-style (s="2", numFmt "YYYY-MM-DD HH:MM:SS") used
# elsewhere in column D to the two brand new rows' date cells.
$ws.Range("D246").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D247").NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Host ("Final dimension: " + $ws.UsedRange.Address())
